$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.445.73"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.852.82"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.18%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.43"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4755"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2756"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06344"
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "17.91"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +11.74%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.845.90"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07465"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.969"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.75"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.94%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.407.49"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "245.53"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +8.14%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.69"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007340"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.926"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.920"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "164.42"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.085"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.878"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1027"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.348"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.046"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.838"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04834"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.130"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.6990"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.699"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01906"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.676"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8793"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.994"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "106.98"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.000"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4064"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.512"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.174"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.27"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +6.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1202"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.04"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.550"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05504"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.347"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3692"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.11%  "
